# Apply the logboek update: fill in row 12 on "week 48" sheet
# (begintijd/eindtijd/activiteit for the 6th entry) and move the
# active-cell selection from F12:F13 to F12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week 48")

# Row 12: begintijd (C12), eindtijd (D12) and activiteit (F12).
$ws.Range("C12").Value = 0.4236111111111111
$ws.Range("D12").Value = 0.42708333333333331
$ws.Range("F12").Value = "Kijken of de Class werkt"

# Update the sheet's active selection to F12 only.
$ws.Activate()
$ws.Range("F12").Select()
